$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14-117 down to 15-118
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new weekly price entry
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44971
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112030
$ws.Range("G14").Value = "Poroto granado"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 35000
$ws.Range("L14").Value = 36000
$ws.Range("M14").Value = 35500
$ws.Range("N14").Value = "`$/malla 25 kilos"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 1420
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
